$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bump the "Dt. Referencia" date column (G) from 2024-11-04 (45600) to
# 2024-11-05 (45601) for every data row (2 through 274).
$ws.Range("G2:G274").Value = 45601

# Row 108: Saldo Previsto / Vl. Total corrected from -136.49 to 24.85
$ws.Range("E108").Value = 24.85
$ws.Range("H108").Value = 24.85

# Row 218: Saldo Previsto / Vl. Total corrected from 0 to 3939.56
$ws.Range("E218").Value = 3939.56
$ws.Range("H218").Value = 3939.56

# Rename the sheet to reflect the new extraction timestamp
$ws.Name = "IClientBalance-20241105-103136-"
